# Apply the "Penalty Reward System" edits to the PO data workbook.
#
# Summary of the change:
#  - On the "Weekly Quantity" sheet, the week of 2023-07-30 (row 11,
#    A11 = 45137.99999999999, B11 = 380) is removed entirely, shifting
#    all subsequent rows up by one.
#  - The requested quantity for the week of 2023-07-23 (now row 10)
#    changes from 290 to 220.
#  - On the "Monthly Trend" sheet, the aggregated requested quantity
#    for the month of 2023-07 (row 4) changes from 1240 to 790 to stay
#    consistent with the weekly-sheet edits above.

$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$monthly = $wb.Worksheets.Item("Monthly Trend")

# Remove the entire row for the week of 2023-07-30 (row 11); this shifts
# every row below it up by one, matching the diff's renumbering.
$weekly.Rows.Item(11).Delete()

# Update the requested quantity for the week of 2023-07-23 (row 10).
$weekly.Cells.Item(10, 2).Value = 220

# Update the Monthly Trend aggregate for the affected month (row 4).
$monthly.Cells.Item(4, 2).Value = 790
